# Workbook / worksheet references
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 241 ("Sin especificar" Espinaca records),
# shifting all of the existing rows 241-277 down to 243-279.
$ws.Range("A241:A242").EntireRow.Insert()

# ---- Row 241 (new record) ----
$ws.Cells.Item(241, 1).Value  = 9
$ws.Cells.Item(241, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(241, 3).Value  = "Metropolitana"
$ws.Cells.Item(241, 4).Value  = 44505
$ws.Cells.Item(241, 5).Value  = 13
$ws.Cells.Item(241, 6).Value  = 100112012
$ws.Cells.Item(241, 7).Value  = "Espinaca"
$ws.Cells.Item(241, 8).Value  = "Sin especificar"
$ws.Cells.Item(241, 9).Value  = "Primera"
$ws.Cells.Item(241, 10).Value = 250
$ws.Cells.Item(241, 11).Value = 5000
$ws.Cells.Item(241, 12).Value = 6000
$ws.Cells.Item(241, 13).Value = 5500
$ws.Cells.Item(241, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(241, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(241, 16).Value = 550
$ws.Cells.Item(241, 17).Value = 10
$ws.Cells.Item(241, 18).Value = "Hortaliza"

# ---- Row 242 (new record) ----
$ws.Cells.Item(242, 1).Value  = 9
$ws.Cells.Item(242, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(242, 3).Value  = "Metropolitana"
$ws.Cells.Item(242, 4).Value  = 44505
$ws.Cells.Item(242, 5).Value  = 13
$ws.Cells.Item(242, 6).Value  = 100112012
$ws.Cells.Item(242, 7).Value  = "Espinaca"
$ws.Cells.Item(242, 8).Value  = "Sin especificar"
$ws.Cells.Item(242, 9).Value  = "Segunda"
$ws.Cells.Item(242, 10).Value = 106
$ws.Cells.Item(242, 11).Value = 4000
$ws.Cells.Item(242, 12).Value = 4000
$ws.Cells.Item(242, 13).Value = 4000
$ws.Cells.Item(242, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(242, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(242, 16).Value = 400
$ws.Cells.Item(242, 17).Value = 10
$ws.Cells.Item(242, 18).Value = "Hortaliza"
